# Add a new worksheet "Criteria Display" with a generic test for
# displaying a single CTS result criterion (ZIP code search).

$wb = $excel.ActiveWorkbook

# Append the new sheet after the last existing sheet (AdvanceSearch) so it
# becomes the last tab and the active tab, matching Excel's default
# behaviour when inserting a new sheet "after" the current last sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Criteria Display"

# Fill cells in the same order the original author entered them, so the
# shared-string table is built up identically.
$ws.Range("C1").Value = "expected text"
$ws.Range("D1").Value = "result URL"
$ws.Range("D2").Value = "https://www-qa.cancer.gov/about-cancer/treatment/clinical-trials/search/r?q=&t=&a=&z=20850&rl=1"
$ws.Range("C2").Value = "within 100 miles of 20850"
$ws.Range("B2").Value = "Near ZIP Code:"
$ws.Range("A2").Value = "Check criteria display for a ZIP code search"
$ws.Range("B1").Value = "expected label"
$ws.Range("A1").Value = "description"

# Approximate the column widths from the original author's sheet
# (closest achievable width given the host's pixel-width quantisation).
$ws.Columns.Item(1).ColumnWidth = 27.67
$ws.Columns.Item(2).ColumnWidth = 22.67
$ws.Columns.Item(3).ColumnWidth = 19
$ws.Columns.Item(4).ColumnWidth = 17.33

# The author's selection ends up on A2 after filling in the data.
$ws.Range("A2").Select() | Out-Null
